# Update the LR-pairs data with newly computed TPM-based values.
# Only the "source" Ligand (G/H) and Receptor (M/N) expression values for the
# "ECs" cluster changed; every other touched column (I, J, O, P, Q, R, S, T)
# is a derived/specificity value that is recomputed from those sources.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New raw expression values (TPM) per cluster.
$G = @{ "ECs" = 5.685592333333333; "FAPs" = 29.800487;          "MuSCs" = 12.35338333333333 }
$H = @{ "ECs" = 17.056777;         "FAPs" = 89.40146100000001;  "MuSCs" = 37.06015 }
$M = @{ "ECs" = 3.626135;          "FAPs" = 0.121294;           "MuSCs" = 3.920705666666667 }
$N = @{ "ECs" = 10.878405;         "FAPs" = 0.363882;           "MuSCs" = 11.762117 }

# Sums across clusters used for the specificity (share-of-total) columns.
$Gsum = 0.0
$Hsum = 0.0
$Msum = 0.0
$Nsum = 0.0
foreach ($k in $G.Keys) { $Gsum += $G[$k] }
foreach ($k in $H.Keys) { $Hsum += $H[$k] }
foreach ($k in $M.Keys) { $Msum += $M[$k] }
foreach ($k in $N.Keys) { $Nsum += $N[$k] }

for ($r = 2; $r -le 10; $r++) {
    $sender = $ws.Cells.Item($r, 1).Value2   # column A - Sending cluster
    $target = $ws.Cells.Item($r, 4).Value2   # column D - Target cluster

    $gVal = $G[$sender]
    $hVal = $H[$sender]
    $mVal = $M[$target]
    $nVal = $N[$target]

    $iVal = $gVal / $Gsum
    $jVal = $hVal / $Hsum
    $oVal = $mVal / $Msum
    $pVal = $nVal / $Nsum

    $qVal = $gVal * $mVal
    $rVal = $hVal * $nVal
    $sVal = $iVal * $oVal
    $tVal = $jVal * $pVal

    $ws.Cells.Item($r, 7).Value  = $gVal   # G
    $ws.Cells.Item($r, 8).Value  = $hVal   # H
    $ws.Cells.Item($r, 9).Value  = $iVal   # I
    $ws.Cells.Item($r, 10).Value = $jVal   # J
    $ws.Cells.Item($r, 13).Value = $mVal   # M
    $ws.Cells.Item($r, 14).Value = $nVal   # N
    $ws.Cells.Item($r, 15).Value = $oVal   # O
    $ws.Cells.Item($r, 16).Value = $pVal   # P
    $ws.Cells.Item($r, 17).Value = $qVal   # Q
    $ws.Cells.Item($r, 18).Value = $rVal   # R
    $ws.Cells.Item($r, 19).Value = $sVal   # S
    $ws.Cells.Item($r, 20).Value = $tVal   # T
}
